# Generate Report for Handoff
#
# - Status text "Handed back: in sync with en-US" -> "Ready for handoff"
#   on the Overview sheet (E2, F2) and on each language sheet's Status column (C2).
# - Refresh the "Latest HO Xliff Generate Date" / "Latest Handback DateTime" and
#   "Latest Handoff Datetime" timestamps to the new handoff run time.
# - Narrow the (now shorter) status columns to their new auto-fit width.

$wb = $excel.ActiveWorkbook

$ws_overview = $wb.Worksheets.Item("Overview")
$ws_zhcn     = $wb.Worksheets.Item("zh-cn")
$ws_dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet -------------------------------------------------------
$ws_overview.Range("E2").Value2 = "Ready for handoff"
$ws_overview.Range("F2").Value2 = "Ready for handoff"
$ws_overview.Range("G2").Value2 = "2016-09-06 11:19:31"

# --- zh-cn sheet ------------------------------------------------------------
$ws_zhcn.Range("C2").Value2 = "Ready for handoff"
$ws_zhcn.Range("H2").Value2 = "2016-09-06 11:19:26"

# --- de-de sheet ------------------------------------------------------------
$ws_dede.Range("C2").Value2 = "Ready for handoff"
$ws_dede.Range("H2").Value2 = "2016-09-06 11:19:31"

# --- Column width adjustments (Status columns shrank after the text change) -
$ws_overview.Columns.Item(5).ColumnWidth = 16.27   # column E
$ws_overview.Columns.Item(6).ColumnWidth = 16.27   # column F
$ws_zhcn.Columns.Item(3).ColumnWidth = 16.27        # column C
$ws_dede.Columns.Item(3).ColumnWidth = 16.27        # column C
